# Apply hybrid bold + color ("#2C3E50") highlighting to quantitative
# metrics (percentages, dollar amounts, large numbers) inside specific
# bullet/impact paragraphs, per the commit:
#   "Implement quantitative metrics highlighting across all resume formats"
#
# Strategy: for each target paragraph, locate it by a stable substring,
# then walk forward through the paragraph, Find()-ing each metric token
# in turn (always re-scoping the search range to [last match end, end of
# paragraph] so repeated tokens like "87%" ... "71%" are not confused),
# and set Font.Bold / Font.Color on just that matched sub-range. Word's
# OM automatically splits the run(s) as needed, mirroring the diff.
#
# NOTE: this PowerShell engine does not bind *named* function arguments
# (e.g. `-Foo bar`) -- only positional. So every call below is positional.

$d = $word.ActiveDocument

# Decimal value that iron_docx/Word renders back out as w:color val="2C3E50".
$HighlightColor = 5258796

function Set-MetricHighlight {
    param(
        [string]$ParaContains,
        [string[]]$Metrics,
        [bool]$ExactMatch
    )

    $paras = $d.Paragraphs
    $target = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $ptext = $p.Range.Text.Trim()
        if ($ExactMatch) {
            if ($ptext -eq $ParaContains) {
                $target = $p
                break
            }
        } elseif ($p.Range.Text -like "*$ParaContains*") {
            $target = $p
            break
        }
    }
    if ($null -eq $target) {
        Write-Output "WARN: paragraph containing '$ParaContains' not found"
        return
    }

    $paraEnd = $target.Range.End
    $cursor = $target.Range.Start

    foreach ($metric in $Metrics) {
        $rng = $d.Range($cursor, $paraEnd)
        $found = $rng.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = 1
            $rng.Font.Color = $HighlightColor
            $cursor = $rng.End
        } else {
            Write-Output "WARN: metric '$metric' not found in paragraph '$ParaContains'"
        }
    }
}

$PM = [char]0xB1   # '±'

Set-MetricHighlight "improving demographic classification accuracy from" @("23%", "64%") $false

Set-MetricHighlight "reducing polling error margins" @("87%", "71%", "$($PM)4.2%", "$($PM)2.1%") $false

Set-MetricHighlight "Wrote RFP and analyzed bids from" @("1,200") $false

Set-MetricHighlight "became the " @("`$400M", "`$1B") $false

Set-MetricHighlight "Algorithm reduced mapping costs by" @("73.5%", "`$4.7M") $false

Set-MetricHighlight "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" @("87%", "71%") $true
